$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing data (A1:A369) shifts down to A2:A370
$ws.Range("A1").EntireRow.Insert()
$ws.Range("A1").Value = "ENSEMBL_ID"
